$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 766.6667
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 766.6667
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 766.6667
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1418.6667
$ws.Range("H88").Value = 2621.25
$ws.Range("I88").Value = 2328.3333
$ws.Range("K88").Value = 2328.3333
$ws.Range("M88").Value = -1922.3333
$ws.Range("H91").Value = 2621.25
$ws.Range("I91").Value = 2328.3333
$ws.Range("K91").Value = 2328.3333
$ws.Range("M91").Value = -924.3332999999998

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 30000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H61").Value = 2699.6
$ws.Range("I61").Value = 2699.6
$ws.Range("K61").Value = 2699.6
$ws.Range("M61").Value = -2487.6
$ws.Range("H74").Value = 7410.4443
$ws.Range("I74").Value = 7167.5
$ws.Range("K74").Value = 7167.5
$ws.Range("M74").Value = -6293.5
$ws.Range("H77").Value = 7410.4443
$ws.Range("I77").Value = 7167.5
$ws.Range("K77").Value = 35837.5
$ws.Range("M77").Value = -31469.5
$ws.Range("H97").Value = 777.5
$ws.Range("I97").Value = 777.5
$ws.Range("K97").Value = 777.5
$ws.Range("M97").Value = -281.5
$ws.Range("H110").Value = 3004.25
$ws.Range("I110").Value = 2855.6667
$ws.Range("K110").Value = 2855.6667
$ws.Range("M110").Value = -810.6667000000002
$ws.Range("H122").Value = 1374.25
$ws.Range("I122").Value = 1149.1666
$ws.Range("K122").Value = 3447.4998
$ws.Range("M122").Value = -997.4998000000001
$ws.Range("H132").Value = 1639.8572
$ws.Range("I132").Value = 1579.8334
$ws.Range("K132").Value = 4739.5002
$ws.Range("M132").Value = -2209.5002
$ws.Range("H136").Value = 2699.6
$ws.Range("I136").Value = 2699.6
$ws.Range("K136").Value = 8098.799999999999
$ws.Range("M136").Value = -5548.799999999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1271
$ws.Range("I99").Value = 1291
$ws.Range("K99").Value = 1291
$ws.Range("M99").Value = 207
$ws.Range("H105").Value = 1054758.1
$ws.Range("I105").Value = 1668867.1
$ws.Range("K105").Value = 1668867.1
$ws.Range("M105").Value = -1667120.1
$ws.Range("H107").Value = 4770.524
$ws.Range("I107").Value = 1409.2727
$ws.Range("J107").Value = 8467.9
$ws.Range("K107").Value = 1409.2727
$ws.Range("L107").Value = 8467.9
$ws.Range("M107").Value = 510.7273
$ws.Range("N107").Value = -12307.9
$ws.Range("H134").Value = 5823.077
$ws.Range("I134").Value = 1190.909
$ws.Range("K134").Value = 3572.727
$ws.Range("M134").Value = -1037.727

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4980.2666
$ws.Range("I31").Value = 2037.2727
$ws.Range("K31").Value = 2037.2727
$ws.Range("M31").Value = -1742.2727
$ws.Range("H34").Value = 4980.2666
$ws.Range("I34").Value = 2037.2727
$ws.Range("K34").Value = 2037.2727
$ws.Range("M34").Value = -1835.2727
$ws.Range("H58").Value = 7829.8335
$ws.Range("J58").Value = 7996
$ws.Range("L58").Value = 7996
$ws.Range("N58").Value = -8402
$ws.Range("H132").Value = 3691.1
$ws.Range("I132").Value = 3691.1
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11073.3
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8543.299999999999
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 3359.6155
$ws.Range("I134").Value = 3232.25
$ws.Range("K134").Value = 9696.75
$ws.Range("M134").Value = -7161.75
$ws.Range("H136").Value = 7829.8335
$ws.Range("J136").Value = 7996
$ws.Range("L136").Value = 23988
$ws.Range("N136").Value = -29088

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 219.33333
$ws.Range("I6").Value = 96.75
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 290.25
$ws.Range("L6").Value = 3600
$ws.Range("M6").Value = -177.25
$ws.Range("N6").Value = -3826
$ws.Range("H55").Value = 1498
$ws.Range("J55").Value = 2032
$ws.Range("L55").Value = 6096
$ws.Range("N55").Value = -6450
$ws.Range("H68").Value = 780
$ws.Range("I68").Value = 1200
$ws.Range("K68").Value = 3600
$ws.Range("M68").Value = -2789
$ws.Range("H71").Value = 780
$ws.Range("I71").Value = 1200
$ws.Range("K71").Value = 10800
$ws.Range("M71").Value = -6744

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 40000
$ws.Range("J15").Value = 40000
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40576
$ws.Range("H75").Value = 74000
$ws.Range("J75").Value = 74000
$ws.Range("L75").Value = 74000
$ws.Range("N75").Value = -75748
$ws.Range("H78").Value = 74000
$ws.Range("J78").Value = 74000
$ws.Range("L78").Value = 222000
$ws.Range("N78").Value = -230736
$ws.Range("H80").Value = 2145.4443
$ws.Range("I80").Value = 2272.8572
$ws.Range("J80").Value = 1699.5
$ws.Range("K80").Value = 2272.8572
$ws.Range("L80").Value = 1699.5
$ws.Range("M80").Value = -1274.8572
$ws.Range("N80").Value = -3695.5
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H83").Value = 2145.4443
$ws.Range("I83").Value = 2272.8572
$ws.Range("J83").Value = 1699.5
$ws.Range("K83").Value = 11364.286
$ws.Range("L83").Value = 8497.5
$ws.Range("M83").Value = -6372.286
$ws.Range("N83").Value = -18481.5
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H102").Value = 1566.1765
$ws.Range("I102").Value = 829.5833
$ws.Range("K102").Value = 829.5833
$ws.Range("M102").Value = 792.4167

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1581.25
$ws.Range("I16").Value = 1559.1052
$ws.Range("K16").Value = 1559.1052
$ws.Range("M16").Value = -1389.1052
$ws.Range("H61").Value = 5435.727
$ws.Range("I61").Value = 4049
$ws.Range("K61").Value = 4049
$ws.Range("M61").Value = -3847
$ws.Range("H82").Value = 3019.8
$ws.Range("I82").Value = 949.8
$ws.Range("J82").Value = 5089.8
$ws.Range("K82").Value = 949.8
$ws.Range("L82").Value = 5089.8
$ws.Range("M82").Value = -588.8
$ws.Range("N82").Value = -5811.8
$ws.Range("H85").Value = 3019.8
$ws.Range("I85").Value = 949.8
$ws.Range("J85").Value = 5089.8
$ws.Range("K85").Value = 949.8
$ws.Range("L85").Value = 5089.8
$ws.Range("M85").Value = 298.2
$ws.Range("N85").Value = -7585.8
$ws.Range("H113").Value = 5435.727
$ws.Range("I113").Value = 4049
$ws.Range("K113").Value = 4049
$ws.Range("M113").Value = -1879
$ws.Range("H132").Value = 3130.3845
$ws.Range("I132").Value = 2886.25
$ws.Range("K132").Value = 8658.75
$ws.Range("M132").Value = -6128.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H122").Value = 2230.5925
$ws.Range("I122").Value = 2063.9565
$ws.Range("K122").Value = 6191.869499999999
$ws.Range("M122").Value = -3741.869499999999
